$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 1).Value = 'data/output/10088/4277181'
$ws.Cells.Item(4, 2).Value = 'N° de subpreguntas incorrecto para estudiante 4277181,                    se encontraron 166 subpreguntas'

$ws.Cells.Item(5, 1).Value = 'data/output/10121/4278294'
$ws.Cells.Item(5, 2).Value = 'N° de subpreguntas incorrecto para estudiante 4278294,                    se encontraron 166 subpreguntas'

$ws.Cells.Item(6, 1).Value = 'data/output/10157/4279607_p20'
$ws.Cells.Item(6, 2).Value = 'Pregunta no pudo ser procesada'

$ws.Cells.Item(7, 1).Value = 'data/output/10157/4279607_p4'
$ws.Cells.Item(7, 2).Value = 'Pregunta no pudo ser procesada'

$ws.Cells.Item(8, 1).Value = 'data/output/10157/4279607'
$ws.Cells.Item(8, 2).Value = 'N° de subpreguntas incorrecto para estudiante 4279607,                    se encontraron 155 subpreguntas'

$ws.Range("A6").Select()

